$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lrpap1"
$ws.Cells.Item(2,3).Value = "Sorl1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 8.191447666666667
$ws.Cells.Item(2,8).Value = 24.574343
$ws.Cells.Item(2,9).Value = 0.185794284429433
$ws.Cells.Item(2,10).Value = 0.185794284429433
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.05271666666666667
$ws.Cells.Item(2,14).Value = 0.15815
$ws.Cells.Item(2,15).Value = 0.003013014833311122
$ws.Cells.Item(2,16).Value = 0.003013014833311122
$ws.Cells.Item(2,17).Value = 0.4318258161611112
$ws.Cells.Item(2,18).Value = 3.88643234545
$ws.Cells.Item(2,19).Value = 0.0005598009349303072
$ws.Cells.Item(2,20).Value = 0.0005598009349303073

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lrpap1"
$ws.Cells.Item(3,3).Value = "Sorl1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 8.191447666666667
$ws.Cells.Item(3,8).Value = 24.574343
$ws.Cells.Item(3,9).Value = 0.185794284429433
$ws.Cells.Item(3,10).Value = 0.185794284429433
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.046984666666667
$ws.Cells.Item(3,14).Value = 3.140954
$ws.Cells.Item(3,15).Value = 0.05984028449413786
$ws.Cells.Item(3,16).Value = 0.05984028449413786
$ws.Cells.Item(3,17).Value = 8.576320104802445
$ws.Cells.Item(3,18).Value = 77.186880943222
$ws.Cells.Item(3,19).Value = 0.01111798283764204
$ws.Cells.Item(3,20).Value = 0.01111798283764204

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lrpap1"
$ws.Cells.Item(4,3).Value = "Sorl1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 8.191447666666667
$ws.Cells.Item(4,8).Value = 24.574343
$ws.Cells.Item(4,9).Value = 0.185794284429433
$ws.Cells.Item(4,10).Value = 0.185794284429433
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 15.25299333333333
$ws.Cells.Item(4,14).Value = 45.75898
$ws.Cells.Item(4,15).Value = 0.8717830255908123
$ws.Cells.Item(4,16).Value = 0.8717830255908123
$ws.Cells.Item(4,17).Value = 124.9440966500156
$ws.Cells.Item(4,18).Value = 1124.49686985014
$ws.Cells.Item(4,19).Value = 0.161972303417371
$ws.Cells.Item(4,20).Value = 0.161972303417371

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Lrpap1"
$ws.Cells.Item(5,3).Value = "Sorl1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 8.191447666666667
$ws.Cells.Item(5,8).Value = 24.574343
$ws.Cells.Item(5,9).Value = 0.185794284429433
$ws.Cells.Item(5,10).Value = 0.185794284429433
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.143623666666667
$ws.Cells.Item(5,14).Value = 3.430871
$ws.Cells.Item(5,15).Value = 0.06536367508173863
$ws.Cells.Item(5,16).Value = 0.06536367508173863
$ws.Cells.Item(5,17).Value = 9.367933415861446
$ws.Cells.Item(5,18).Value = 84.31140074275301
$ws.Cells.Item(5,19).Value = 0.01214419723948959
$ws.Cells.Item(5,20).Value = 0.01214419723948959

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lrpap1"
$ws.Cells.Item(6,3).Value = "Sorl1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.317702
$ws.Cells.Item(6,8).Value = 51.95310600000001
$ws.Cells.Item(6,9).Value = 0.3927913821808575
$ws.Cells.Item(6,10).Value = 0.3927913821808576
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.05271666666666667
$ws.Cells.Item(6,14).Value = 0.15815
$ws.Cells.Item(6,15).Value = 0.003013014833311122
$ws.Cells.Item(6,16).Value = 0.003013014833311122
$ws.Cells.Item(6,17).Value = 0.9129315237666668
$ws.Cells.Item(6,18).Value = 8.2163837139
$ws.Cells.Item(6,19).Value = 0.001183486260907702
$ws.Cells.Item(6,20).Value = 0.001183486260907702

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lrpap1"
$ws.Cells.Item(7,3).Value = "Sorl1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.317702
$ws.Cells.Item(7,8).Value = 51.95310600000001
$ws.Cells.Item(7,9).Value = 0.3927913821808575
$ws.Cells.Item(7,10).Value = 0.3927913821808576
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.046984666666667
$ws.Cells.Item(7,14).Value = 3.140954
$ws.Cells.Item(7,15).Value = 0.05984028449413786
$ws.Cells.Item(7,16).Value = 0.05984028449413786
$ws.Cells.Item(7,17).Value = 18.13136845590267
$ws.Cells.Item(7,18).Value = 163.182316103124
$ws.Cells.Item(7,19).Value = 0.02350474805654814
$ws.Cells.Item(7,20).Value = 0.02350474805654815

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Lrpap1"
$ws.Cells.Item(8,3).Value = "Sorl1"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 17.317702
$ws.Cells.Item(8,8).Value = 51.95310600000001
$ws.Cells.Item(8,9).Value = 0.3927913821808575
$ws.Cells.Item(8,10).Value = 0.3927913821808576
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 15.25299333333333
$ws.Cells.Item(8,14).Value = 45.75898
$ws.Cells.Item(8,15).Value = 0.8717830255908123
$ws.Cells.Item(8,16).Value = 0.8717830255908123
$ws.Cells.Item(8,17).Value = 264.1467931546534
$ws.Cells.Item(8,18).Value = 2377.32113839188
$ws.Cells.Item(8,19).Value = 0.3424288595836251
$ws.Cells.Item(8,20).Value = 0.3424288595836251

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Lrpap1"
$ws.Cells.Item(9,3).Value = "Sorl1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 17.317702
$ws.Cells.Item(9,8).Value = 51.95310600000001
$ws.Cells.Item(9,9).Value = 0.3927913821808575
$ws.Cells.Item(9,10).Value = 0.3927913821808576
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.143623666666667
$ws.Cells.Item(9,14).Value = 3.430871
$ws.Cells.Item(9,15).Value = 0.06536367508173863
$ws.Cells.Item(9,16).Value = 0.06536367508173863
$ws.Cells.Item(9,17).Value = 19.80493385948067
$ws.Cells.Item(9,18).Value = 178.244404735326
$ws.Cells.Item(9,19).Value = 0.02567428827977659
$ws.Cells.Item(9,20).Value = 0.02567428827977659

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Lrpap1"
$ws.Cells.Item(10,3).Value = "Sorl1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 12.35128266666667
$ws.Cells.Item(10,8).Value = 37.053848
$ws.Cells.Item(10,9).Value = 0.2801455637905346
$ws.Cells.Item(10,10).Value = 0.2801455637905346
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.05271666666666667
$ws.Cells.Item(10,14).Value = 0.15815
$ws.Cells.Item(10,15).Value = 0.003013014833311122
$ws.Cells.Item(10,16).Value = 0.003013014833311122
$ws.Cells.Item(10,17).Value = 0.6511184512444446
$ws.Cells.Item(10,18).Value = 5.8600660612
$ws.Cells.Item(10,19).Value = 0.0008440827391871878
$ws.Cells.Item(10,20).Value = 0.0008440827391871879

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Lrpap1"
$ws.Cells.Item(11,3).Value = "Sorl1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 12.35128266666667
$ws.Cells.Item(11,8).Value = 37.053848
$ws.Cells.Item(11,9).Value = 0.2801455637905346
$ws.Cells.Item(11,10).Value = 0.2801455637905346
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.046984666666667
$ws.Cells.Item(11,14).Value = 3.140954
$ws.Cells.Item(11,15).Value = 0.05984028449413786
$ws.Cells.Item(11,16).Value = 0.05984028449413786
$ws.Cells.Item(11,17).Value = 12.93160356566578
$ws.Cells.Item(11,18).Value = 116.384432090992
$ws.Cells.Item(11,19).Value = 0.01676399023699623
$ws.Cells.Item(11,20).Value = 0.01676399023699623

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Lrpap1"
$ws.Cells.Item(12,3).Value = "Sorl1"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 12.35128266666667
$ws.Cells.Item(12,8).Value = 37.053848
$ws.Cells.Item(12,9).Value = 0.2801455637905346
$ws.Cells.Item(12,10).Value = 0.2801455637905346
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 15.25299333333333
$ws.Cells.Item(12,14).Value = 45.75898
$ws.Cells.Item(12,15).Value = 0.8717830255908123
$ws.Cells.Item(12,16).Value = 0.8717830255908123
$ws.Cells.Item(12,17).Value = 188.3940321727823
$ws.Cells.Item(12,18).Value = 1695.54628955504
$ws.Cells.Item(12,19).Value = 0.2442261472071562
$ws.Cells.Item(12,20).Value = 0.2442261472071562

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Lrpap1"
$ws.Cells.Item(13,3).Value = "Sorl1"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 12.35128266666667
$ws.Cells.Item(13,8).Value = 37.053848
$ws.Cells.Item(13,9).Value = 0.2801455637905346
$ws.Cells.Item(13,10).Value = 0.2801455637905346
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.143623666666667
$ws.Cells.Item(13,14).Value = 3.430871
$ws.Cells.Item(13,15).Value = 0.06536367508173863
$ws.Cells.Item(13,16).Value = 0.06536367508173863
$ws.Cells.Item(13,17).Value = 14.12521917128978
$ws.Cells.Item(13,18).Value = 127.126972541608
$ws.Cells.Item(13,19).Value = 0.01831134360719498
$ws.Cells.Item(13,20).Value = 0.01831134360719498

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Lrpap1"
$ws.Cells.Item(14,3).Value = "Sorl1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 6.228371
$ws.Cells.Item(14,8).Value = 18.685113
$ws.Cells.Item(14,9).Value = 0.1412687695991749
$ws.Cells.Item(14,10).Value = 0.1412687695991749
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.05271666666666667
$ws.Cells.Item(14,14).Value = 0.15815
$ws.Cells.Item(14,15).Value = 0.003013014833311122
$ws.Cells.Item(14,16).Value = 0.003013014833311122
$ws.Cells.Item(14,17).Value = 0.3283389578833333
$ws.Cells.Item(14,18).Value = 2.95505062095
$ws.Cells.Item(14,19).Value = 0.0004256448982859252
$ws.Cells.Item(14,20).Value = 0.0004256448982859252

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Lrpap1"
$ws.Cells.Item(15,3).Value = "Sorl1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 6.228371
$ws.Cells.Item(15,8).Value = 18.685113
$ws.Cells.Item(15,9).Value = 0.1412687695991749
$ws.Cells.Item(15,10).Value = 0.1412687695991749
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 1.046984666666667
$ws.Cells.Item(15,14).Value = 3.140954
$ws.Cells.Item(15,15).Value = 0.05984028449413786
$ws.Cells.Item(15,16).Value = 0.05984028449413786
$ws.Cells.Item(15,17).Value = 6.521008935311333
$ws.Cells.Item(15,18).Value = 58.689080417802
$ws.Cells.Item(15,19).Value = 0.008453563362951437
$ws.Cells.Item(15,20).Value = 0.008453563362951437

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Lrpap1"
$ws.Cells.Item(16,3).Value = "Sorl1"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 6.228371
$ws.Cells.Item(16,8).Value = 18.685113
$ws.Cells.Item(16,9).Value = 0.1412687695991749
$ws.Cells.Item(16,10).Value = 0.1412687695991749
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 15.25299333333333
$ws.Cells.Item(16,14).Value = 45.75898
$ws.Cells.Item(16,15).Value = 0.8717830255908123
$ws.Cells.Item(16,16).Value = 0.8717830255908123
$ws.Cells.Item(16,17).Value = 95.00130134052668
$ws.Cells.Item(16,18).Value = 855.0117120647401
$ws.Cells.Item(16,19).Value = 0.12315571538266
$ws.Cells.Item(16,20).Value = 0.12315571538266

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Lrpap1"
$ws.Cells.Item(17,3).Value = "Sorl1"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 6.228371
$ws.Cells.Item(17,8).Value = 18.685113
$ws.Cells.Item(17,9).Value = 0.1412687695991749
$ws.Cells.Item(17,10).Value = 0.1412687695991749
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 1.143623666666667
$ws.Cells.Item(17,14).Value = 3.430871
$ws.Cells.Item(17,15).Value = 0.06536367508173863
$ws.Cells.Item(17,16).Value = 0.06536367508173863
$ws.Cells.Item(17,17).Value = 7.122912480380334
$ws.Cells.Item(17,18).Value = 64.10621232342301
$ws.Cells.Item(17,19).Value = 0.00923384595527746
$ws.Cells.Item(17,20).Value = 0.00923384595527746
